$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 edits
$ws.Range("B2").Value = "Test 5th"
$ws.Range("C2").Value = "tablets"
$ws.Range("E2").Value = "mg"
$ws.Range("F2").Value = "Test 5th"
$ws.Range("G2").Value = "Test 5th"
$ws.Range("H2").Value = "64745"
$ws.Range("I2").Value = "2017-03-10"
$ws.Range("J2").Value = "EMS"
$ws.Range("K2").Value = 12
$ws.Range("L2").Value = 24
$ws.Range("M2").Value = 12
$ws.Range("N2").Value = "21"

# Remove bold-ish style (s4) from B2,C2,F2,G2 -- handled implicitly by not copying style

# New Row 3
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Test 20th"
$ws.Range("C3").Value = "tablets"
$ws.Range("D3").Value = 1000
$ws.Range("E3").Value = "mg"
$ws.Range("F3").Value = "Test 20th"
$ws.Range("G3").Value = "Test 20th"
$ws.Range("H3").Value = "64745sd"
$ws.Range("I3").Value = "2017-03-10"
$ws.Range("J3").Value = "EMS"
$ws.Range("K3").Value = 12
$ws.Range("L3").Value = 24
$ws.Range("M3").Value = 12
$ws.Range("N3").Value = "21"

# Apply style 4 (same as row 6 style) to row 3
$ws.Range("A3:N3").Style = $ws.Range("A6:N6").Style
